$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Slit3"
$ws.Range("C2").Value = "Robo2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8417399999999999
$ws.Range("H2").Value = 1.68348
$ws.Range("I2").Value = 0.0119744574213572
$ws.Range("J2").Value = 0.008423523143208028
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5215465
$ws.Range("N2").Value = 1.043093
$ws.Range("O2").Value = 0.3993854683363039
$ws.Range("P2").Value = 0.3071469942766215
$ws.Range("Q2").Value = 0.43900655091
$ws.Range("R2").Value = 1.75602620364
$ws.Range("S2").Value = 0.004782424285301875
$ws.Range("T2").Value = 0.002587259814655905

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Slit3"
$ws.Range("C3").Value = "Robo2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8417399999999999
$ws.Range("H3").Value = 1.68348
$ws.Range("I3").Value = 0.0119744574213572
$ws.Range("J3").Value = 0.008423523143208028
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.002342
$ws.Range("N3").Value = 0.007026
$ws.Range("O3").Value = 0.001793436954986034
$ws.Range("P3").Value = 0.002068861340060323
$ws.Range("Q3").Value = 0.00197135508
$ws.Range("R3").Value = 0.01182813048
$ws.Range("S3").Value = 0.00002147543445536877
$ws.Range("T3").Value = 0.0000174271013780865

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Slit3"
$ws.Range("C4").Value = "Robo2"
$ws.Range("D4").Value = "Neutro"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8417399999999999
$ws.Range("H4").Value = 1.68348
$ws.Range("I4").Value = 0.0119744574213572
$ws.Range("J4").Value = 0.008423523143208028
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.7819839999999999
$ws.Range("N4").Value = 2.345952
$ws.Range("O4").Value = 0.59882109470871
$ws.Range("P4").Value = 0.6907841443833183
$ws.Range("Q4").Value = 0.6582272121599999
$ws.Range("R4").Value = 3.949363272959999
$ws.Range("S4").Value = 0.007170557701599954
$ws.Range("T4").Value = 0.005818836227174037

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Slit3"
$ws.Range("C5").Value = "Robo2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 57.09845533333333
$ws.Range("H5").Value = 171.295366
$ws.Range("I5").Value = 0.8122734124721006
$ws.Range("J5").Value = 0.8570998644624763
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5215465
$ws.Range("N5").Value = 1.043093
$ws.Range("O5").Value = 0.3993854683363039
$ws.Range("P5").Value = 0.3071469942766215
$ws.Range("Q5").Value = 29.77949953450634
$ws.Range("R5").Value = 178.676997207038
$ws.Range("S5").Value = 0.3244101972572976
$ws.Range("T5").Value = 0.2632556471645492

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Slit3"
$ws.Range("C6").Value = "Robo2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 57.09845533333333
$ws.Range("H6").Value = 171.295366
$ws.Range("I6").Value = 0.8122734124721006
$ws.Range("J6").Value = 0.8570998644624763
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.002342
$ws.Range("N6").Value = 0.007026
$ws.Range("O6").Value = 0.001793436954986034
$ws.Range("P6").Value = 0.002068861340060323
$ws.Range("Q6").Value = 0.1337245823906667
$ws.Range("R6").Value = 1.203521241516
$ws.Range("S6").Value = 0.001456761155480079
$ws.Range("T6").Value = 0.00177322077415736

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Slit3"
$ws.Range("C7").Value = "Robo2"
$ws.Range("D7").Value = "Neutro"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 57.09845533333333
$ws.Range("H7").Value = 171.295366
$ws.Range("I7").Value = 0.8122734124721006
$ws.Range("J7").Value = 0.8570998644624763
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7819839999999999
$ws.Range("N7").Value = 2.345952
$ws.Range("O7").Value = 0.59882109470871
$ws.Range("P7").Value = 0.6907841443833183
$ws.Range("Q7").Value = 44.65007849538133
$ws.Range("R7").Value = 401.8507064584319
$ws.Range("S7").Value = 0.4864064540593229
$ws.Range("T7").Value = 0.5920709965237697

# Row 8
$ws.Range("A8").Value = "M1"
$ws.Range("B8").Value = "Slit3"
$ws.Range("C8").Value = "Robo2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.02523366666666666
$ws.Range("H8").Value = 0.07570099999999999
$ws.Range("I8").Value = 0.000358970070431155
$ws.Range("J8").Value = 0.0003787803392163797
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5215465
$ws.Range("N8").Value = 1.043093
$ws.Range("O8").Value = 0.3993854683363039
$ws.Range("P8").Value = 0.3071469942766215
$ws.Range("Q8").Value = 0.01316053053216667
$ws.Range("R8").Value = 0.078963183193
$ws.Range("S8").Value = 0.0001433674296978628
$ws.Range("T8").Value = 0.0001163412426813901

# Row 9
$ws.Range("A9").Value = "M1"
$ws.Range("B9").Value = "Slit3"
$ws.Range("C9").Value = "Robo2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.02523366666666666
$ws.Range("H9").Value = 0.07570099999999999
$ws.Range("I9").Value = 0.000358970070431155
$ws.Range("J9").Value = 0.0003787803392163797
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.002342
$ws.Range("N9").Value = 0.007026
$ws.Range("O9").Value = 0.001793436954986034
$ws.Range("P9").Value = 0.002068861340060323
$ws.Range("Q9").Value = 0.00005909724733333333
$ws.Range("R9").Value = 0.000531875226
$ws.Range("S9").Value = 0.0000006437901900451729
$ws.Range("T9").Value = 0.000000783644000179703

# Row 10
$ws.Range("A10").Value = "M1"
$ws.Range("B10").Value = "Slit3"
$ws.Range("C10").Value = "Robo2"
$ws.Range("D10").Value = "Neutro"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.02523366666666666
$ws.Range("H10").Value = 0.07570099999999999
$ws.Range("I10").Value = 0.000358970070431155
$ws.Range("J10").Value = 0.0003787803392163797
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.7819839999999999
$ws.Range("N10").Value = 2.345952
$ws.Range("O10").Value = 0.59882109470871
$ws.Range("P10").Value = 0.6907841443833183
$ws.Range("Q10").Value = 0.01973232359466666
$ws.Range("R10").Value = 0.1775909123519999
$ws.Range("S10").Value = 0.000214958850543247
$ws.Range("T10").Value = 0.0002616554525348099

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Slit3"
$ws.Range("C11").Value = "Robo2"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.08110000000000001
$ws.Range("H11").Value = 0.2433
$ws.Range("I11").Value = 0.00115371551413984
$ws.Range("J11").Value = 0.001217384929278942
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.5215465
$ws.Range("N11").Value = 1.043093
$ws.Range("O11").Value = 0.3993854683363039
$ws.Range("P11").Value = 0.3071469942766215
$ws.Range("Q11").Value = 0.04229742115000001
$ws.Range("R11").Value = 0.2537845269
$ws.Range("S11").Value = 0.0004607772109415997
$ws.Range("T11").Value = 0.0003739161219056845

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Slit3"
$ws.Range("C12").Value = "Robo2"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.08110000000000001
$ws.Range("H12").Value = 0.2433
$ws.Range("I12").Value = 0.00115371551413984
$ws.Range("J12").Value = 0.001217384929278942
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.002342
$ws.Range("N12").Value = 0.007026
$ws.Range("O12").Value = 0.001793436954986034
$ws.Range("P12").Value = 0.002068861340060323
$ws.Range("Q12").Value = 0.0001899362
$ws.Range("R12").Value = 0.0017094258
$ws.Range("S12").Value = 0.000002069116038599101
$ws.Range("T12").Value = 0.000002518600616157274

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Slit3"
$ws.Range("C13").Value = "Robo2"
$ws.Range("D13").Value = "Neutro"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.08110000000000001
$ws.Range("H13").Value = 0.2433
$ws.Range("I13").Value = 0.00115371551413984
$ws.Range("J13").Value = 0.001217384929278942
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.7819839999999999
$ws.Range("N13").Value = 2.345952
$ws.Range("O13").Value = 0.59882109470871
$ws.Range("P13").Value = 0.6907841443833183
$ws.Range("Q13").Value = 0.06341890239999999
$ws.Range("R13").Value = 0.5707701216
$ws.Range("S13").Value = 0.0006908691871596413
$ws.Range("T13").Value = 0.0008409502067571005

# Row 14
$ws.Range("A14").Value = "Neutro"
$ws.Range("B14").Value = "Slit3"
$ws.Range("C14").Value = "Robo2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 2.060578
$ws.Range("H14").Value = 6.181734
$ws.Range("I14").Value = 0.02931345014420769
$ws.Range("J14").Value = 0.03093115416527428
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.5215465
$ws.Range("N14").Value = 1.043093
$ws.Range("O14").Value = 0.3993854683363039
$ws.Range("P14").Value = 0.3071469942766215
$ws.Range("Q14").Value = 1.074687243877
$ws.Range("R14").Value = 6.448123463262
$ws.Range("S14").Value = 0.01170736601439728
$ws.Range("T14").Value = 0.009500411031370793

# Row 15
$ws.Range("A15").Value = "Neutro"
$ws.Range("B15").Value = "Slit3"
$ws.Range("C15").Value = "Robo2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2.060578
$ws.Range("H15").Value = 6.181734
$ws.Range("I15").Value = 0.02931345014420769
$ws.Range("J15").Value = 0.03093115416527428
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.002342
$ws.Range("N15").Value = 0.007026
$ws.Range("O15").Value = 0.001793436954986034
$ws.Range("P15").Value = 0.002068861340060323
$ws.Range("Q15").Value = 0.004825873676
$ws.Range("R15").Value = 0.043432863084
$ws.Range("S15").Value = 0.00005257182476676275
$ws.Range("T15").Value = 0.00006399226905598177

# Row 16
$ws.Range("A16").Value = "Neutro"
$ws.Range("B16").Value = "Slit3"
$ws.Range("C16").Value = "Robo2"
$ws.Range("D16").Value = "Neutro"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 2.060578
$ws.Range("H16").Value = 6.181734
$ws.Range("I16").Value = 0.02931345014420769
$ws.Range("J16").Value = 0.03093115416527428
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.7819839999999999
$ws.Range("N16").Value = 2.345952
$ws.Range("O16").Value = 0.59882109470871
$ws.Range("P16").Value = 0.6907841443833183
$ws.Range("Q16").Value = 1.611339026752
$ws.Range("R16").Value = 14.502051240768
$ws.Range("S16").Value = 0.01755351230504364
$ws.Range("T16").Value = 0.0213667508648475

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Slit3"
$ws.Range("C17").Value = "Robo2"
$ws.Range("D17").Value = "ECs"
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 10.1875185
$ws.Range("H17").Value = 20.375037
$ws.Range("I17").Value = 0.1449259943777636
$ws.Range("J17").Value = 0.1019492929605459
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.5215465
$ws.Range("N17").Value = 1.043093
$ws.Range("O17").Value = 0.3993854683363039
$ws.Range("P17").Value = 0.3071469942766215
$ws.Range("Q17").Value = 5.31326461736025
$ws.Range("R17").Value = 21.253058469441
$ws.Range("S17").Value = 0.05788133613866767
$ws.Range("T17").Value = 0.0313134189014584

# Row 18
$ws.Range("A18").Value = "sCs"
$ws.Range("B18").Value = "Slit3"
$ws.Range("C18").Value = "Robo2"
$ws.Range("D18").Value = "FAPs"
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 10.1875185
$ws.Range("H18").Value = 20.375037
$ws.Range("I18").Value = 0.1449259943777636
$ws.Range("J18").Value = 0.1019492929605459
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = 0.3333333333333333
$ws.Range("M18").Value = 0.002342
$ws.Range("N18").Value = 0.007026
$ws.Range("O18").Value = 0.001793436954986034
$ws.Range("P18").Value = 0.002068861340060323
$ws.Range("Q18").Value = 0.023859168327
$ws.Range("R18").Value = 0.143155009962
$ws.Range("S18").Value = 0.0002599156340551795
$ws.Range("T18").Value = 0.0002109189508525575

# Row 19
$ws.Range("A19").Value = "sCs"
$ws.Range("B19").Value = "Slit3"
$ws.Range("C19").Value = "Robo2"
$ws.Range("D19").Value = "Neutro"
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 10.1875185
$ws.Range("H19").Value = 20.375037
$ws.Range("I19").Value = 0.1449259943777636
$ws.Range("J19").Value = 0.1019492929605459
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 0.7819839999999999
$ws.Range("N19").Value = 2.345952
$ws.Range("O19").Value = 0.59882109470871
$ws.Range("P19").Value = 0.6907841443833183
$ws.Range("Q19").Value = 7.966476466703998
$ws.Range("R19").Value = 47.79885880022399
$ws.Range("S19").Value = 0.08678474260504077
$ws.Range("T19").Value = 0.07042495510823497
